$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename/extend columns (Name/Floors -> building_name..postcode) ---
$ws.Range("A1").Value = "building_name"
$ws.Range("B1").Value = "street_address1"
$ws.Range("C1").Value = "street_address2"
$ws.Range("D1").Value = "city"
$ws.Range("E1").Value = "county"
$ws.Range("F1").Value = "postcode"

# Header cells keep their original style slot, but are no longer bold.
$ws.Range("A1:B1").Font.Bold = $false

# --- Row 2: Chaucer House ---
$ws.Range("A2").Value = "Chaucer House"
$ws.Range("D2").Value = "Portsmouth"
$ws.Range("E2").Value = "Hampshire"

# --- Row 3: Catherine House ---
$ws.Range("A3").Value = "Catherine House"
$ws.Range("B3").Value = "Stanhope Rd"
$ws.Range("D3").Value = "Portsmouth"
$ws.Range("E3").Value = "Hampshire"
$ws.Range("F3").Value = "PO11DZ"

# back to finish off row 2
$ws.Range("F2").Value = "PO12DR"

# --- Row 4: Greetham Street Hall ---
$ws.Range("A4").Value = "Greetham Street Hall"
$ws.Range("B4").Value = "Greetham St"
$ws.Range("D4").Value = "Portsmouth"
$ws.Range("E4").Value = "Hampshire"
$ws.Range("F4").Value = "PO54FB"

# --- Row 5: Margaret Rule Hall ---
$ws.Range("A5").Value = "Margaret Rule Hall"
$ws.Range("D5").Value = "Portsmouth"
$ws.Range("E5").Value = "Hampshire"
$ws.Range("F5").Value = "PO12DS"
$ws.Range("B5").Value = "Isambard Brunel Road"

# back to finish off row 2
$ws.Range("B2").Value = "32 Isambard Brunel Road"

# --- Row 6: Harry Law Hall ---
$ws.Range("A6").Value = "Harry Law Hall"
$ws.Range("B6").Value = "Isambard Brunel Road"
$ws.Range("D6").Value = "Portsmouth"
$ws.Range("E6").Value = "Hampshire"
$ws.Range("F6").Value = "PO12SP"

# --- Row 7: James Watson Hall ---
$ws.Range("A7").Value = "James Watson Hall"
$ws.Range("B7").Value = "28 Guildhall Walk"
$ws.Range("D7").Value = "Portsmouth"
$ws.Range("E7").Value = "Hampshire"
$ws.Range("F7").Value = "PO12DD"

# --- Row 8: Trafalgar Hall ---
$ws.Range("A8").Value = "Trafalgar Hall"
$ws.Range("B8").Value = "Middle Street"
$ws.Range("C8").Value = "Southsea"
$ws.Range("D8").Value = "Portsmouth"
$ws.Range("E8").Value = "Hampshire"
$ws.Range("F8").Value = "PO54AY"

# --- Row 9: Bateson Hall ---
$ws.Range("A9").Value = "Bateson Hall"
$ws.Range("B9").Value = "The Mary Rose Street"
$ws.Range("D9").Value = "Portsmouth"
$ws.Range("E9").Value = "Hampshire"
$ws.Range("F9").Value = "PO12BL"

# --- Row 10: Rees Hall ---
$ws.Range("A10").Value = "Rees Hall"
$ws.Range("B10").Value = "Southsea Terrace"
$ws.Range("C10").Value = "Southsea"
$ws.Range("D10").Value = "Portsmouth"
$ws.Range("E10").Value = "Hampshire"
$ws.Range("F10").Value = "PO53AP"

# Column widths for the newly-populated columns (best-fit to content).
$ws.Range("B:E").EntireColumn.AutoFit() | Out-Null

# Restore the cursor/selection position as left by the editor.
$ws.Range("B17").Select()
